# Apply the data-format change described in the commit:
#   "modify data format, replace data in 2 diff dfs."
#
# Summary of the change to the 항목설명 (item description) sheet:
#  - Rows 22-24 (bond yield rows) are renamed to explicitly be USA bond yields
#      10_bond  / 10_YEAR_BOND_YIELD  / 10년만기국채수익률   -> bond_usa_10 / 10_YEAR_BOND_YIELD_usa / 미국10년만기국채수익률
#      2_bond   / 2_YEAR_BOND_YIELD   / 2년만기국채수익률    -> bond_usa_2  / 2_YEAR_BOND_YIELD_usa  / 미국2년만기국채수익률
#      3_m_bond / 3_MONTH_BOND_YIELD  / 3개월만기국채수익률  -> bond_usa_3m / 3_MONTH_BOND_YIELD_usa / 미국3개월만기국채수익률
#  - Two new rows are inserted right after those for Korean bond yields:
#      bond_kor_10 / 10_YEAR_BOND_YIELD_korea / 미국10년만기국채수익률
#      bond_kor_2  / 2_YEAR_BOND_YIELD_korea  / 미국2년만기국채수익률
#  - Every row below shifts down by 2, and column A (순번/번호) is kept in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new blank rows at 25 and 26 (this pushes the old rows 25.. down to 27..)
$ws.Rows("25:26").Insert()

# 2. Update the (now-USA-specific) bond yield rows 22-24
$ws.Range("B22").Value = "bond_usa_10"
$ws.Range("C22").Value = "10_YEAR_BOND_YIELD_usa"
$ws.Range("D22").Value = "미국10년만기국채수익률"

$ws.Range("B23").Value = "bond_usa_2"
$ws.Range("C23").Value = "2_YEAR_BOND_YIELD_usa"
$ws.Range("D23").Value = "미국2년만기국채수익률"

$ws.Range("B24").Value = "bond_usa_3m"
$ws.Range("C24").Value = "3_MONTH_BOND_YIELD_usa"
$ws.Range("D24").Value = "미국3개월만기국채수익률"

# 3. Fill the two newly inserted rows with the Korean bond yield data
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "bond_kor_10"
$ws.Range("C25").Value = "10_YEAR_BOND_YIELD_korea"
$ws.Range("D25").Value = "미국10년만기국채수익률"

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "bond_kor_2"
$ws.Range("C26").Value = "2_YEAR_BOND_YIELD_korea"
$ws.Range("D26").Value = "미국2년만기국채수익률"

# 4. Re-sync the running index (column A) for every data row so it always equals row-1.
for ($r = 2; $r -le 67; $r++) {
    $ws.Range("A$r").Value = $r - 1
}

# 5. Restore the sheet view state to match the saved workbook (scrolled down a bit,
#    selection sitting on the newly appended row).
$ws.Range("B68").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
